$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Merge the three runs "We also " / "create a matrix" / " to map out..."
#    into a single run. Doing a Find/Replace with identical replacement
#    text that starts inside the "create a matrix" run merges the runs
#    it spans without altering the visible text.
# ------------------------------------------------------------------
$d.Content.Find.Execute("create a matrix", $true, $false, $false, $false, $false, $true, 1, $false, "create a matrix", 2) | Out-Null

# ------------------------------------------------------------------
# 2. "Shared Movements/Actions" heading: font size 24 -> 32 (half-points),
#    i.e. Font.Size 12 -> 16 (points). Paragraph index 3 (1-based) in
#    Word's Paragraphs collection.
# ------------------------------------------------------------------
$d.Paragraphs(3).Range.Font.Size = 16

# ------------------------------------------------------------------
# 3. "Defender Movements" heading: font size 24 -> 32.
# ------------------------------------------------------------------
$d.Paragraphs(10).Range.Font.Size = 16

# ------------------------------------------------------------------
# 4. "The defender will also go into Hunting Mode if: " ->
#    "The defender will also go into Hunting Mode if" | " one of the
#    following is true" | (bookmark _GoBack) | ": "
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("The defender will also go into Hunting Mode if: ") | Out-Null
$huntStart = $rng.Start
$huntEnd = $rng.End

# Remove the trailing ": " first -- a direct Range.Text="" delete
# correctly recomputes xml:space on the shrunk run.
$tailRange = $d.Range($huntEnd - 2, $huntEnd)
$tailRange.Text = ""
$ifEnd = $huntEnd - 2

# Insert the replacement tail text.
$insertPoint = $d.Range($ifEnd, $ifEnd)
$insertPoint.InsertAfter(" one of the following is true: ")

# Split out the "The defender ... if" run from " one of ... true: ".
$splitBk = $d.Range($ifEnd, $ifEnd)
$d.Bookmarks.Add("TempSplitHunt", $splitBk) | Out-Null
$d.Bookmarks("TempSplitHunt").Delete()

# Split out ": " into its own run, marked by the relocated _GoBack bookmark.
$trueEnd = $ifEnd + 29
$gobackRange = $d.Range($trueEnd, $trueEnd)
$d.Bookmarks.Add("_GoBack", $gobackRange) | Out-Null

Write-Host "stage2 ok"
